# Commit: "Added script for gcash OTP and MPIN"
#
# Net effect observed in the target OOXML:
#  - test_suite: PayThruCash runmode N->Y, BaseTest runmode Y->N,
#    and the extra "HappyPath" run-control row is removed (13 rows -> 12).
#  - OrderDetails: trimmed back down to a single "Quantity" column
#    (the Payment Method / E-Wallet / Gcash helper columns are removed).
#  - the active sheet/selection moves off payThruGcash and back onto
#    test_suite, and several other sheets get their remembered selection
#    bumped around too (harmless view state that tags along with normal
#    interactive editing).
#  - payThruMasterCard's column D width loses its "best fit" flag.

$wb = $excel.ActiveWorkbook

# ---- test_suite: flip the two runmode flags ---------------------------
$ts = $wb.Worksheets.Item("test_suite")
$ts.Range("B2").Value = "Y"    # PayThruCash
$ts.Range("B12").Value = "N"   # BaseTest

# ---- test_suite: drop the trailing HappyPath row -----------------------
$ts.Rows.Item(13).Delete()

# ---- OrderDetails: remove the Payment Method / E-Wallet helper columns -
$od = $wb.Worksheets.Item("OrderDetails")
$od.Range("B1:C2").ClearContents()

# ---- payThruMasterCard: column D no longer best-fit ---------------------
$mc = $wb.Worksheets.Item("payThruMasterCard")
$mc.Columns.Item(4).ColumnWidth = 5.5

# ---- view/selection bookkeeping ----------------------------------------
$od.Range("E11").Select()

$gc = $wb.Worksheets.Item("payThruGcash")
$gc.Select()
$gc.Range("C2").Select()

$vi = $wb.Worksheets.Item("payThruVisa")
$vi.Select()
$vi.Range("E20").Select()

$mc.Select()
$mc.Range("G18").Select()

# test_suite becomes the active sheet again, zoomed to 140% with G3 selected
$ts.Select()
$ts.Application.ActiveWindow.Zoom = 140
$ts.Range("G3").Select()
